$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the table header columns / cells (Table1 headers live in A1:B1)
$ws.Range("A1").Value = "MaterialID"
$ws.Range("B1").Value = "Quantity"

# Update the active selection to match the saved workbook view
$ws.Range("F7").Select()
